$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7098885070416259
$ws.Range("C2").Value = 0.07864104582836262
$ws.Range("D2").Value = 0.1565908986480764
$ws.Range("F2").Value = 1.696331049648762
$ws.Range("G2").Value = 0.002490942961933952
$ws.Range("I2").Value = 1.06570494582818
$ws.Range("J2").Value = 0.1914952365474001
$ws.Range("K2").Value = 0.4151878403228864
$ws.Range("L2").Value = 0.3667765988833622
$ws.Range("N2").Value = 1.986254401642821
$ws.Range("O2").Value = 4.290002999546601

$ws.Range("B3").Value = 0.6723398164128582
$ws.Range("C3").Value = 0.07687536464467826
$ws.Range("D3").Value = 0.1544067579605723
$ws.Range("F3").Value = 1.702367124285772
$ws.Range("G3").Value = 0.00249323044903841
$ws.Range("I3").Value = 1.072873905163434
$ws.Range("J3").Value = 0.1923844011265778
$ws.Range("K3").Value = 0.3805412918438265
$ws.Range("L3").Value = 0.3616697045864186
$ws.Range("N3").Value = 2.00431188527564
$ws.Range("O3").Value = 4.312658234789382

$ws.Range("B4").Value = 0.6495078773059788
$ws.Range("C4").Value = 0.0757781368124455
$ws.Range("D4").Value = 0.153121906508531
$ws.Range("F4").Value = 1.706796981319542
$ws.Range("G4").Value = 0.002494711434941876
$ws.Range("I4").Value = 1.077701877552371
$ws.Range("J4").Value = 0.1929921663213463
$ws.Range("K4").Value = 0.3593497618215338
$ws.Range("L4").Value = 0.3586810303730488
$ws.Range("N4").Value = 2.015969442555933
$ws.Range("O4").Value = 4.328422246415798

$ws.Range("B5").Value = 0.640260550521532
$ws.Range("C5").Value = 0.07532773068972176
$ws.Range("D5").Value = 0.1526125373765623
$ws.Range("F5").Value = 1.70878439786955
$ws.Range("G5").Value = 0.002495334230740354
$ws.Range("I5").Value = 1.079776584002772
$ws.Range("J5").Value = 0.1932554075663937
$ws.Range("K5").Value = 0.3507351545163431
$ws.Range("L5").Value = 0.3575002435648571
$ws.Range("N5").Value = 2.020863431661582
$ws.Range("O5").Value = 4.335312740473086

$ws.Range("B6").Value = 0.638728496680983
$ws.Range("C6").Value = 0.07525274363059253
$ws.Range("D6").Value = 0.1525288181601354
$ws.Range("F6").Value = 1.709125419140264
$ws.Range("G6").Value = 0.002495438811862291
$ws.Range("I6").Value = 1.08012756902211
$ws.Range("J6").Value = 0.1933000599148009
$ws.Range("K6").Value = 0.349305998277174
$ws.Range("L6").Value = 0.3573064216967623
$ws.Range("N6").Value = 2.021684740072608
$ws.Range("O6").Value = 4.336485091686143

$ws.Range("B7").Value = 0.6493829333203109
$ws.Range("C7").Value = 0.07577207572475686
$ws.Range("D7").Value = 0.1531149793014279
$ws.Range("F7").Value = 1.706823046221395
$ws.Range("G7").Value = 0.002494719756043959
$ws.Range("I7").Value = 1.077729423315031
$ws.Range("J7").Value = 0.1929956533965758
$ws.Range("K7").Value = 0.3592334959373886
$ws.Range("L7").Value = 0.3586649553458727
$ws.Range("N7").Value = 2.016034863739929
$ws.Range("O7").Value = 4.328513284494477

$ws.Range("B8").Value = 0.6968958623337187
$ws.Range("C8").Value = 0.07803496502012308
$ws.Range("D8").Value = 0.1558261863672854
$ws.Range("F8").Value = 1.698262228507211
$ws.Range("G8").Value = 0.002491715851983011
$ws.Range("I8").Value = 1.068088402379143
$ws.Range("J8").Value = 0.1917890077348225
$ws.Range("K8").Value = 0.4032251231651571
$ws.Range("L8").Value = 0.3649853316194793
$ws.Range("N8").Value = 1.992362317538773
$ws.Range("O8").Value = 4.297430061855465

$ws.Range("B9").Value = 0.7918099514203334
$ws.Range("C9").Value = 0.08236811154957735
$ws.Range("D9").Value = 0.1615857444580655
$ws.Range("F9").Value = 1.687207101441018
$ws.Range("G9").Value = 0.002486429336515674
$ws.Range("I9").Value = 1.052560021999643
$ws.Range("J9").Value = 0.1899120462737258
$ws.Range("K9").Value = 0.4901180793824267
$ws.Range("L9").Value = 0.378539566834462
$ws.Range("N9").Value = 1.950461308353457
$ws.Range("O9").Value = 4.251166382742412

$ws.Range("B10").Value = 0.8625734577283311
$ws.Range("C10").Value = 0.08548765853898743
$ws.Range("D10").Value = 0.1660835684304232
$ws.Range("F10").Value = 1.682568323281416
$ws.Range("G10").Value = 0.002482910070171381
$ws.Range("I10").Value = 1.043205129972762
$ws.Range("J10").Value = 0.1888297480673842
$ws.Range("K10").Value = 0.5543166396610388
$ws.Range("L10").Value = 0.3891980550330771
$ws.Range("N10").Value = 1.922426958398121
$ws.Range("O10").Value = 4.226111604660844

$ws.Range("B11").Value = 0.8949830498382596
$ws.Range("C11").Value = 0.08689285430975957
$ws.Range("D11").Value = 0.1681868089715408
$ws.Range("F11").Value = 1.68121211944954
$ws.Range("G11").Value = 0.002481387507693591
$ws.Range("I11").Value = 1.039394222359434
$ws.Range("J11").Value = 0.1884014831653396
$ws.Range("K11").Value = 0.5835955177280994
$ws.Range("L11").Value = 0.3941975203448322
$ws.Range("N11").Value = 1.910269407163038
$ws.Range("O11").Value = 4.21664958291916

$ws.Range("B12").Value = 0.9072865085987019
$ws.Range("C12").Value = 0.0874229541388658
$ws.Range("D12").Value = 0.1689913952262003
$ws.Range("F12").Value = 1.68080677268766
$ws.Range("G12").Value = 0.002480822165392724
$ws.Range("I12").Value = 1.038014987942077
$ws.Range("J12").Value = 0.1882484976141754
$ws.Range("K12").Value = 0.5946928724337113
$ws.Range("L12").Value = 0.3961122348914188
$ws.Range("N12").Value = 1.905751222301125
$ws.Range("O12").Value = 4.213344519954575

$ws.Range("B13").Value = 0.904635388696164
$ws.Range("C13").Value = 0.0873088776680504
$ws.Range("D13").Value = 0.168817752733176
$ws.Range("F13").Value = 1.680889261614155
$ws.Range("G13").Value = 0.00248094342373947
$ws.Range("I13").Value = 1.038309191081346
$ws.Range("J13").Value = 0.1882810375051314
$ws.Range("K13").Value = 0.5923024206733771
$ws.Range("L13").Value = 0.3956989121906105
$ws.Range("N13").Value = 1.906720486496715
$ws.Range("O13").Value = 4.214043965499911

$ws.Range("B14").Value = 0.8959946532899608
$ws.Range("C14").Value = 0.08693650647700935
$ws.Range("D14").Value = 0.1682528402193384
$ws.Range("F14").Value = 1.681176603654109
$ws.Range("G14").Value = 0.002481340772146167
$ws.Range("I14").Value = 1.03927947219664
$ws.Range("J14").Value = 0.1883887129339392
$ws.Range("K14").Value = 0.5845083055210409
$ws.Range("L14").Value = 0.3943546146264367
$ws.Range("N14").Value = 1.909895977518105
$ws.Range("O14").Value = 4.216372103748114

$ws.Range("B15").Value = 0.8907059201870595
$ws.Range("C15").Value = 0.08670815511599272
$ws.Range("D15").Value = 0.1679078718637754
$ws.Range("F15").Value = 1.681366696037756
$ws.Range("G15").Value = 0.00248158561907712
$ws.Range("I15").Value = 1.039882113525735
$ws.Range("J15").Value = 0.1884558631832043
$ws.Range("K15").Value = 0.5797354781424531
$ws.Range("L15").Value = 0.3935339918009788
$ws.Range("N15").Value = 1.911852207106408
$ws.Range("O15").Value = 4.217834349683301

$ws.Range("B16").Value = 0.8604597488642014
$ws.Range("C16").Value = 0.08539554439159502
$ws.Range("D16").Value = 0.1659472596524125
$ws.Range("F16").Value = 1.682672109180942
$ws.Range("G16").Value = 0.002483011145903311
$ws.Range("I16").Value = 1.043463123780725
$ws.Range("J16").Value = 0.1888590234374803
$ws.Range("K16").Value = 0.5524046420921422
$ws.Range("L16").Value = 0.3888743495981828
$ws.Range("N16").Value = 1.923233463833472
$ws.Range("O16").Value = 4.226768889779095

$ws.Range("B17").Value = 0.841960181748874
$ws.Range("C17").Value = 0.08458672665902611
$ws.Range("D17").Value = 0.1647590669330441
$ws.Range("F17").Value = 1.683665908680915
$ws.Range("G17").Value = 0.00248390569674517
$ws.Range("I17").Value = 1.0457738002493
$ws.Range("J17").Value = 0.1891227431034679
$ws.Range("K17").Value = 0.5356567098455685
$ws.Range("L17").Value = 0.3860543303239439
$ws.Range("N17").Value = 1.930367992598273
$ws.Range("O17").Value = 4.232745442070808

$ws.Range("B18").Value = 0.8313403753164152
$ws.Range("C18").Value = 0.08412020868008341
$ws.Range("D18").Value = 0.1640810349318684
$ws.Range("F18").Value = 1.684308498201041
$ws.Range("G18").Value = 0.002484427598061885
$ws.Range("I18").Value = 1.047144699809987
$ws.Range("J18").Value = 0.1892804605096288
$ws.Range("K18").Value = 0.5260308058562941
$ws.Range("L18").Value = 0.3844465401089678
$ws.Range("N18").Value = 1.934527628206887
$ws.Range("O18").Value = 4.236365196893416

$ws.Range("B19").Value = 0.8277482641997551
$ws.Range("C19").Value = 0.0839620297149466
$ws.Range("D19").Value = 0.1638523924032143
$ws.Range("F19").Value = 1.68453826553344
$ws.Range("G19").Value = 0.002484605574137893
$ws.Range("I19").Value = 1.047616054430801
$ws.Range("J19").Value = 0.1893348978284521
$ws.Range("K19").Value = 0.5227728732396599
$ws.Range("L19").Value = 0.3839046157359576
$ws.Range("N19").Value = 1.935945635021239
$ws.Range("O19").Value = 4.237622086667955

$ws.Range("B20").Value = 0.8439273605214623
$ws.Range("C20").Value = 0.08467296213181896
$ws.Range("D20").Value = 0.1648849953588467
$ws.Range("F20").Value = 1.68355277211446
$ws.Range("G20").Value = 0.002483809706778004
$ws.Range("I20").Value = 1.045523492999436
$ws.Range("J20").Value = 0.1890940455032784
$ws.Range("K20").Value = 0.5374388294881101
$ws.Range("L20").Value = 0.3863530567011964
$ws.Range("N20").Value = 1.929602709590222
$ws.Range("O20").Value = 4.232090373222405

$ws.Range("B21").Value = 0.8985318215325719
$ws.Range("C21").Value = 0.08704593581686026
$ws.Range("D21").Value = 0.1684185486186323
$ws.Range("F21").Value = 1.681089268963433
$ws.Range("G21").Value = 0.002481223757216076
$ws.Range("I21").Value = 1.038992744165334
$ws.Range("J21").Value = 0.1883568368625319
$ws.Range("K21").Value = 0.5867973577368275
$ws.Range("L21").Value = 0.3947488845565346
$ws.Range("N21").Value = 1.908960935257454
$ws.Range("O21").Value = 4.215680730888351

$ws.Range("B22").Value = 0.934397257121276
$ws.Range("C22").Value = 0.08858504437361603
$ws.Range("D22").Value = 0.1707753059955337
$ws.Range("F22").Value = 1.680109910033522
$ws.Range("G22").Value = 0.002479599059174147
$ws.Range("I22").Value = 1.035096793758349
$ws.Range("J22").Value = 0.1879285776750343
$ws.Range("K22").Value = 0.6191144542147526
$ws.Range("L22").Value = 0.4003614322311932
$ws.Range("N22").Value = 1.895969377647522
$ws.Range("O22").Value = 4.206576322537018

$ws.Range("B23").Value = 0.9152391625122789
$ws.Range("C23").Value = 0.08776467587794201
$ws.Range("D23").Value = 0.1695131532519554
$ws.Range("F23").Value = 1.680574972905404
$ws.Range("G23").Value = 0.002480460226158383
$ws.Range("I23").Value = 1.037142096272667
$ws.Range("J23").Value = 0.1881522562286904
$ws.Range("K23").Value = 0.6018610797174517
$ws.Range("L23").Value = 0.397354492636083
$ws.Range("N23").Value = 1.902857555352361
$ws.Range("O23").Value = 4.211287372183506

$ws.Range("B24").Value = 0.8430379495918316
$ws.Range("C24").Value = 0.08463397980138865
$ws.Range("D24").Value = 0.1648280473016968
$ws.Range("F24").Value = 1.683603699230808
$ws.Range("G24").Value = 0.002483853080114273
$ws.Range("I24").Value = 1.04563652463905
$ws.Range("J24").Value = 0.1891070006815099
$ws.Range("K24").Value = 0.5366331247217317
$ws.Range("L24").Value = 0.3862179605168734
$ws.Range("N24").Value = 1.929948513681948
$ws.Range("O24").Value = 4.232385957438453

$ws.Range("B25").Value = 0.7659500013709248
$ws.Range("C25").Value = 0.08120709443146268
$ws.Range("D25").Value = 0.1599806015982921
$ws.Range("F25").Value = 1.689585262566695
$ws.Range("G25").Value = 0.002487795171887663
$ws.Range("I25").Value = 1.056399809506246
$ws.Range("J25").Value = 0.1903675967726883
$ws.Range("K25").Value = 0.4665466461734127
$ws.Range("L25").Value = 0.3747493217197189
$ws.Range("N25").Value = 1.911852207106408
$ws.Range("O25").Value = 4.217834349683301
